$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the first five MP records (rows 2-6) with a newer data pull.
# Columns: A=Name, B=Name Url, C=Political Affiliation, D=Constituency,
#          E=Constituency Url, F=Province, G=Phone, H=Email

$ws.Range("A2").Value = "Anstey, Carol"
$ws.Range("B2").Value = "https://www.ourcommons.ca/Members/en/carol-anstey(109872)"
$ws.Range("C2").Value = "Conservative"
$ws.Range("D2").Value = "Long Range Mountains"
$ws.Range("E2").Value = "https://www.ourcommons.ca/Members/en/constituencies/long-range-mountains(947)"
$ws.Range("F2").Value = "Newfoundland and Labrador"
$ws.Range("G2").Value = "709-637-4655"
$ws.Range("H2").Value = "carol.anstey@parl.gc.ca"

$ws.Range("A3").Value = "Arnold, Mel"
$ws.Range("B3").Value = "https://www.ourcommons.ca/Members/en/mel-arnold(89294)"
$ws.Range("C3").Value = "Conservative"
$ws.Range("D3").Value = "Kamloops—Shuswap—Central Rockies"
$ws.Range("E3").Value = "https://www.ourcommons.ca/Members/en/constituencies/kamloops-shuswap-central-rockies(1253)"
$ws.Range("F3").Value = "British Columbia"
$ws.Range("G3").Value = "778-283-9700"
$ws.Range("H3").Value = "mel.arnold@parl.gc.ca"

$ws.Range("A4").Value = "Au, Chak"
$ws.Range("B4").Value = "https://www.ourcommons.ca/Members/en/chak-au(123608)"
$ws.Range("C4").Value = "Conservative"
$ws.Range("D4").Value = "Richmond Centre—Marpole"
$ws.Range("E4").Value = "https://www.ourcommons.ca/Members/en/constituencies/richmond-centre-marpole(1266)"
$ws.Range("F4").Value = "British Columbia"
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = "chak.au@parl.gc.ca"

$ws.Range("A5").Value = "Auguste, Tatiana"
$ws.Range("B5").Value = "https://www.ourcommons.ca/Members/en/tatiana-auguste(122753)"
$ws.Range("C5").Value = "Liberal"
$ws.Range("D5").Value = "Terrebonne"
$ws.Range("E5").Value = "https://www.ourcommons.ca/Members/en/constituencies/terrebonne(1047)"
$ws.Range("F5").Value = "Quebec"
$ws.Range("G5").Value = "450-964-4919"
$ws.Range("H5").Value = "tatiana.auguste@parl.gc.ca"

$ws.Range("A6").Value = "Baber, Roman"
$ws.Range("B6").Value = "https://www.ourcommons.ca/Members/en/roman-baber(123276)"
$ws.Range("C6").Value = "Conservative"
$ws.Range("D6").Value = "York Centre"
$ws.Range("E6").Value = "https://www.ourcommons.ca/Members/en/constituencies/york-centre(1172)"
$ws.Range("F6").Value = "Ontario"
$ws.Range("G6").Value = "416-638-9030"
$ws.Range("H6").Value = "roman.baber@parl.gc.ca"

$wb.Save()
